$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.318.98"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = "  -1.87%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.007.18"
$ws.Range("D3").ClearFormats()

$ws.Range("E3").Value = "  -1.97%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.62"
$ws.Range("D5").ClearFormats()

$ws.Range("E5").Value = "  -1.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.65"
$ws.Range("D6").ClearFormats()

$ws.Range("E6").Value = "  -5.87%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -2.48%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.007.00"
$ws.Range("D9").ClearFormats()

$ws.Range("E9").Value = "  -1.97%  "

$ws.Range("E10").Value = "  -5.13%  "

$ws.Range("E11").Value = "  -1.34%  "

$ws.Range("E12").Value = "  +2.39%  "

$ws.Range("E13").Value = "  -3.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.47"
$ws.Range("D14").ClearFormats()

$ws.Range("E14").Value = "  -6.78%  "

$ws.Range("E15").Value = "  +2.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.498.01"
$ws.Range("D16").ClearFormats()

$ws.Range("E16").Value = "  -2.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.10"
$ws.Range("D17").ClearFormats()

$ws.Range("E17").Value = "  -1.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.266.00"
$ws.Range("D18").ClearFormats()

$ws.Range("E18").Value = "  -1.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.007.75"
$ws.Range("D19").ClearFormats()

$ws.Range("E19").Value = "  -2.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "455.43"
$ws.Range("D20").ClearFormats()

$ws.Range("E20").Value = "  -7.39%  "

$ws.Range("E21").Value = "  -3.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.688"
$ws.Range("D22").ClearFormats()

$ws.Range("E22").Value = "  -2.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.41"
$ws.Range("D23").ClearFormats()

$ws.Range("E23").Value = "  -2.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.57"
$ws.Range("D24").ClearFormats()

$ws.Range("E24").Value = "  -0.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.38"
$ws.Range("D25").ClearFormats()

$ws.Range("E25").Value = "  -3.85%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.22"
$ws.Range("D26").ClearFormats()

$ws.Range("E26").Value = "  -9.48%  "

$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.01"
$ws.Range("D28").ClearFormats()

$ws.Range("E28").Value = "  -6.97%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").ClearFormats()

$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("E30").Value = "  -3.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.99"
$ws.Range("D31").ClearFormats()

$ws.Range("E31").Value = "  -5.69%  "

$ws.Range("E32").Value = "  -5.80%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.22"
$ws.Range("D33").ClearFormats()

$ws.Range("E33").Value = "  +3.27%  "

$ws.Range("E34").Value = "  -3.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0799"
$ws.Range("D35").ClearFormats()

$ws.Range("E35").Value = "  -3.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.02"
$ws.Range("D36").ClearFormats()

$ws.Range("E36").Value = "  -3.95%  "

$ws.Range("E37").Value = "  -3.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.12"
$ws.Range("D38").ClearFormats()

$ws.Range("E38").Value = "  -5.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.16"
$ws.Range("D39").ClearFormats()

$ws.Range("E39").Value = "  -1.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.23"
$ws.Range("D40").ClearFormats()

$ws.Range("E40").Value = "  -0.75%  "

$ws.Range("E41").Value = "  -13.15%  "

$ws.Range("E42").Value = "  +3.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "390.92"
$ws.Range("D43").ClearFormats()

$ws.Range("E43").Value = "  -10.81%  "

$ws.Range("E44").Value = "  -1.74%  "

$ws.Range("E45").Value = "  -7.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.720.61"
$ws.Range("D46").ClearFormats()

$ws.Range("E46").Value = "  -4.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "37.14"
$ws.Range("D47").ClearFormats()

$ws.Range("E47").Value = "  -5.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.46"
$ws.Range("D48").ClearFormats()

$ws.Range("E48").Value = "  -1.16%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Stellar"
$ws.Range("B50").ClearFormats()

$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C50").ClearFormats()

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.109"
$ws.Range("D50").ClearFormats()

$ws.Range("E50").Value = "  -0.91%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("B51").ClearFormats()

$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("C51").ClearFormats()

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.20"
$ws.Range("D51").ClearFormats()

$ws.Range("E51").Value = "  -1.78%  "
